$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Widen column A from 64 to 66 characters ---
# (ColumnWidth setter adds a ~0.8333 offset vs the stored OOXML width, so
#  compensate to land exactly on 66.)
$ws.Columns.Item(1).ColumnWidth = 65.16666666666667

# --- 2. Insert a new row above the existing "Bad Drivers" data row (row 3). ---
# This pushes the old row 3 (AE adapter) down to row 4, "Totals:" down to
# row 5, and the whole "Good Drivers" block down by one row as well.
$ws.Rows.Item(3).EntireRow.Insert()

# New row 3: updated Realtek RTL8852BE driver entry
$ws.Range("A3").Value = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.123.330"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 98.2

# Row 4 (previously row 3): update Critical Minutes / Good Roaming values
$ws.Range("C4").Value = 209
$ws.Range("D4").Value = 98.2

# Row 5 (previously row 4, "Totals:"): update the rolled-up totals
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 210

# --- 3. Append the new "Good Drivers" rows (13-20) ---
$driverRows = @(
    @{ Row = 13; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.123.322"; Samples = 16989;   Pct = 99.9; Vintage = "2024-06-30" },
    @{ Row = 14; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.152.0";   Samples = 1033024; Pct = 100;  Vintage = "2024-04-15" },
    @{ Row = 15; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.149.0";   Samples = 81427;   Pct = 100;  Vintage = "2023-12-20" },
    @{ Row = 16; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.144.0";   Samples = 17672;   Pct = 100;  Vintage = "2023-07-10" },
    @{ Row = 17; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.143.0";   Samples = 326032;  Pct = 100;  Vintage = "2023-06-05" },
    @{ Row = 18; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.141.0";   Samples = 48191;   Pct = 100;  Vintage = "2023-04-17" },
    @{ Row = 19; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.128.0";   Samples = 82442;   Pct = 99.9; Vintage = "2022-08-29" },
    @{ Row = 20; Name = "Realtek RTL8852BE WiFi 6 802.11ax PCIe Adapter - 6001.15.124.0";   Samples = 11789;   Pct = 99.9; Vintage = "2022-07-03" }
)

# Use a scratch row, far away from the real data, as a staging area to get
# vintage strings like "2024-06-30" into cells as literal text (instead of
# Excel auto-converting them to dates) and then copy just the *value* of
# that staged/Text-formatted cell over onto the real destination cell
# (whose own number format is left completely untouched).
$scratchRow = 200
$ws.Rows.Item($scratchRow).EntireRow.Insert()
$scratchCell = $ws.Cells.Item($scratchRow, 1)
$scratchCell.NumberFormat = "@"

foreach ($d in $driverRows) {
    $r = $d.Row

    $ws.Cells.Item($r, 1).Value = $d.Name

    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4152
    $ws.Cells.Item($r, 2).Value = $d.Samples

    $ws.Cells.Item($r, 4).HorizontalAlignment = -4152
    $ws.Cells.Item($r, 4).Value = $d.Pct

    $scratchCell.Value = $d.Vintage
    $scratchCell.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4163)
    $ws.Cells.Item($r, 5).HorizontalAlignment = -4152
}

$ws.Rows.Item($scratchRow).EntireRow.Delete()
